$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95

# Column A holds a date string; force it to be treated as text (matching
# the existing rows which store dates as plain strings, not date serials)
# and avoid leaving behind a number-format style.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "02/27/2026"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = 9228.780000000001
$ws.Cells.Item($row, 3).Value = 0.2451365678766332
$ws.Cells.Item($row, 4).Value = 0.7548634321233668
$ws.Cells.Item($row, 5).Value = -342.25
$ws.Cells.Item($row, 6).Value = -34.75
$ws.Cells.Item($row, 7).Value = -24075.29
$ws.Cells.Item($row, 8).Value = -77.56
$ws.Cells.Item($row, 9).Value = -1190.54
$ws.Cells.Item($row, 10).Value = -34.48
$ws.Cells.Item($row, 11).Value = -25265.83
$ws.Cells.Item($row, 12).Value = -73.25
